$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows before row 836, shifting existing data down (old 836-888 -> new 844-896)
$ws.Rows("836:843").Insert()

# Row 836
$ws.Cells.Item(836, 1).Value = 6
$ws.Cells.Item(836, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(836, 3).Value = 'Metropolitana'
$ws.Cells.Item(836, 4).Value = "2022-07-04"
$ws.Cells.Item(836, 5).Value = 13
$ws.Cells.Item(836, 6).Value = 100112021
$ws.Cells.Item(836, 7).Value = 'Ají'
$ws.Cells.Item(836, 8).Value = 'Americana (o)'
$ws.Cells.Item(836, 9).Value = 'Primera'
$ws.Cells.Item(836, 10).Value = 25
$ws.Cells.Item(836, 11).Value = 40000
$ws.Cells.Item(836, 12).Value = 45000
$ws.Cells.Item(836, 13).Value = 43000
$ws.Cells.Item(836, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(836, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(836, 16).Value = 1720
$ws.Cells.Item(836, 17).Value = 25
$ws.Cells.Item(836, 18).Value = 'Hortaliza'

# Row 837
$ws.Cells.Item(837, 1).Value = 6
$ws.Cells.Item(837, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(837, 3).Value = 'Metropolitana'
$ws.Cells.Item(837, 4).Value = "2022-07-04"
$ws.Cells.Item(837, 5).Value = 13
$ws.Cells.Item(837, 6).Value = 100112021
$ws.Cells.Item(837, 7).Value = 'Ají'
$ws.Cells.Item(837, 8).Value = 'Americana (o)'
$ws.Cells.Item(837, 9).Value = 'Primera'
$ws.Cells.Item(837, 10).Value = 40
$ws.Cells.Item(837, 11).Value = 35000
$ws.Cells.Item(837, 12).Value = 40000
$ws.Cells.Item(837, 13).Value = 38125
$ws.Cells.Item(837, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(837, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(837, 16).Value = 1525
$ws.Cells.Item(837, 17).Value = 25
$ws.Cells.Item(837, 18).Value = 'Hortaliza'

# Row 838
$ws.Cells.Item(838, 1).Value = 6
$ws.Cells.Item(838, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(838, 3).Value = 'Metropolitana'
$ws.Cells.Item(838, 4).Value = "2022-07-04"
$ws.Cells.Item(838, 5).Value = 13
$ws.Cells.Item(838, 6).Value = 100112021
$ws.Cells.Item(838, 7).Value = 'Ají'
$ws.Cells.Item(838, 8).Value = 'Americana (o)'
$ws.Cells.Item(838, 9).Value = 'Segunda'
$ws.Cells.Item(838, 10).Value = 10
$ws.Cells.Item(838, 11).Value = 35000
$ws.Cells.Item(838, 12).Value = 35000
$ws.Cells.Item(838, 13).Value = 35000
$ws.Cells.Item(838, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(838, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(838, 16).Value = 1400
$ws.Cells.Item(838, 17).Value = 25
$ws.Cells.Item(838, 18).Value = 'Hortaliza'

# Row 839
$ws.Cells.Item(839, 1).Value = 6
$ws.Cells.Item(839, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(839, 3).Value = 'Metropolitana'
$ws.Cells.Item(839, 4).Value = "2022-07-04"
$ws.Cells.Item(839, 5).Value = 13
$ws.Cells.Item(839, 6).Value = 100112021
$ws.Cells.Item(839, 7).Value = 'Ají'
$ws.Cells.Item(839, 8).Value = 'Americana (o)'
$ws.Cells.Item(839, 9).Value = 'Segunda'
$ws.Cells.Item(839, 10).Value = 15
$ws.Cells.Item(839, 11).Value = 30000
$ws.Cells.Item(839, 12).Value = 30000
$ws.Cells.Item(839, 13).Value = 30000
$ws.Cells.Item(839, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(839, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(839, 16).Value = 1200
$ws.Cells.Item(839, 17).Value = 25
$ws.Cells.Item(839, 18).Value = 'Hortaliza'

# Row 840
$ws.Cells.Item(840, 1).Value = 6
$ws.Cells.Item(840, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(840, 3).Value = 'Metropolitana'
$ws.Cells.Item(840, 4).Value = "2022-07-04"
$ws.Cells.Item(840, 5).Value = 13
$ws.Cells.Item(840, 6).Value = 100112021
$ws.Cells.Item(840, 7).Value = 'Ají'
$ws.Cells.Item(840, 8).Value = 'Inferno'
$ws.Cells.Item(840, 9).Value = 'Primera'
$ws.Cells.Item(840, 10).Value = 65
$ws.Cells.Item(840, 11).Value = 10000
$ws.Cells.Item(840, 12).Value = 12000
$ws.Cells.Item(840, 13).Value = 10923
$ws.Cells.Item(840, 14).Value = '$/caja 12 kilos'
$ws.Cells.Item(840, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(840, 16).Value = 910
$ws.Cells.Item(840, 17).Value = 12
$ws.Cells.Item(840, 18).Value = 'Hortaliza'

# Row 841
$ws.Cells.Item(841, 1).Value = 6
$ws.Cells.Item(841, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(841, 3).Value = 'Metropolitana'
$ws.Cells.Item(841, 4).Value = "2022-07-04"
$ws.Cells.Item(841, 5).Value = 13
$ws.Cells.Item(841, 6).Value = 100112021
$ws.Cells.Item(841, 7).Value = 'Ají'
$ws.Cells.Item(841, 8).Value = 'Inferno'
$ws.Cells.Item(841, 9).Value = 'Primera'
$ws.Cells.Item(841, 10).Value = 28
$ws.Cells.Item(841, 11).Value = 15000
$ws.Cells.Item(841, 12).Value = 17000
$ws.Cells.Item(841, 13).Value = 16071
$ws.Cells.Item(841, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(841, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(841, 16).Value = 1071
$ws.Cells.Item(841, 17).Value = 15
$ws.Cells.Item(841, 18).Value = 'Hortaliza'

# Row 842
$ws.Cells.Item(842, 1).Value = 6
$ws.Cells.Item(842, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(842, 3).Value = 'Metropolitana'
$ws.Cells.Item(842, 4).Value = "2022-07-04"
$ws.Cells.Item(842, 5).Value = 13
$ws.Cells.Item(842, 6).Value = 100112021
$ws.Cells.Item(842, 7).Value = 'Ají'
$ws.Cells.Item(842, 8).Value = 'Inferno'
$ws.Cells.Item(842, 9).Value = 'Segunda'
$ws.Cells.Item(842, 10).Value = 25
$ws.Cells.Item(842, 11).Value = 8000
$ws.Cells.Item(842, 12).Value = 8000
$ws.Cells.Item(842, 13).Value = 8000
$ws.Cells.Item(842, 14).Value = '$/caja 12 kilos'
$ws.Cells.Item(842, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(842, 16).Value = 667
$ws.Cells.Item(842, 17).Value = 12
$ws.Cells.Item(842, 18).Value = 'Hortaliza'

# Row 843
$ws.Cells.Item(843, 1).Value = 6
$ws.Cells.Item(843, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(843, 3).Value = 'Metropolitana'
$ws.Cells.Item(843, 4).Value = "2022-07-04"
$ws.Cells.Item(843, 5).Value = 13
$ws.Cells.Item(843, 6).Value = 100112021
$ws.Cells.Item(843, 7).Value = 'Ají'
$ws.Cells.Item(843, 8).Value = 'Inferno'
$ws.Cells.Item(843, 9).Value = 'Segunda'
$ws.Cells.Item(843, 10).Value = 10
$ws.Cells.Item(843, 11).Value = 12000
$ws.Cells.Item(843, 12).Value = 12000
$ws.Cells.Item(843, 13).Value = 12000
$ws.Cells.Item(843, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(843, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(843, 16).Value = 800
$ws.Cells.Item(843, 17).Value = 15
$ws.Cells.Item(843, 18).Value = 'Hortaliza'
